# Generate Report for Handback
# Update handback/handoff timestamps and the zh-cn/de-de "Priority" (ht -> mt)
# values to reflect a freshly-regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview.Range("G2").Value = "2016-08-25 00:15:09"
$wsOverview.Range("G5").Value = "2016-08-25 00:15:09"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime (column H)
$wsZhCn.Range("H2").Value = "2016-08-25 00:14:58"
$wsZhCn.Range("H5").Value = "2016-08-25 00:14:58"

# Correspond Handback DateTime (column K)
$wsZhCn.Range("K2").Value = "2016-08-25 00:15:29"
$wsZhCn.Range("K5").Value = "2016-08-25 00:15:29"

# --- de-de sheet ---
# Correspond Handoff Datetime (column H) -- shares text with Overview's G column
$wsDeDe.Range("H2").Value = "2016-08-25 00:15:09"
$wsDeDe.Range("H5").Value = "2016-08-25 00:15:09"

# Correspond Handback DateTime (column K)
$wsDeDe.Range("K2").Value = "2016-08-25 00:15:36"
$wsDeDe.Range("K5").Value = "2016-08-25 00:15:36"
